$wb = $excel.ActiveWorkbook

# --- Add the new "ZDB-Nummer" sheet as a copy of "Sprachangaben", placed   ---
# --- at the end of the tab strip (mirrors how the real fixture sheet was   ---
# --- produced: same styles / column widths / layout, new text + data).    ---
$src  = $wb.Worksheets.Item("Sprachangaben")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $last)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "ZDB-Nummer"

# --- Header block (column A labels stay the same, column B gets new text) ---
$ws.Range("B1").Value = "ZDB-Nummer"
$ws.Range("B2").Value = "Nazeige der ZDB-Nummer"
$ws.Range("B3").Value = "System Control Number"
$ws.Range("B4").Value = "ZDB-Nummer"
$ws.Range("B5").Value = "ZDB-Number"
$ws.Range("B6").Value = "Ticket #132"

# --- Replace the Redmine ticket hyperlink (was issues/134, now issues/132) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "http://redmine.thulb.uni-jena.de/issues/132", [type]::Missing, [type]::Missing, "Ticket #132")

# Hyperlinks.Add silently recolors the host cell ("Hyperlink" style); restore
# the original cell formatting (style index 5, shared with B1..B5) by pulling
# formats back in from a still-untouched sibling cell.
$ws.Range("B4").Copy()
$ws.Range("B6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Test-data row: 035 $a / ZDB number example ----------------------------
$ws.Range("A13").Value = "035 `$a"
$ws.Range("B13").Value = 233814418
$ws.Range("C13").Value = "1393051-5"
$ws.Range("D13").ClearContents()

# --- Selection / active-cell bookkeeping on the new sheet ------------------
$ws.Range("C13").Select()

# --- Move the workbook's active tab to the newly appended sheet ------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
